$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("3_3")
$ws2 = $wb.Worksheets.Item("1_8")

# --- Sheet "3_3": rows 2-4, columns B-E get the "Call BT func" measurement values ---
$ws1.Range("B2").Value = 6.8955
$ws1.Range("C2").Value = 6.9009
$ws1.Range("D2").Value = 6.8895
$ws1.Range("E2").Value = 0.0032

$ws1.Range("B3").Value = 6.8957
$ws1.Range("C3").Value = 6.9031
$ws1.Range("D3").Value = 6.892
$ws1.Range("E3").Value = 0.0034

$ws1.Range("B4").Value = 6.8974
$ws1.Range("C4").Value = 6.9027
$ws1.Range("D4").Value = 6.8931
$ws1.Range("E4").Value = 0.003

# Sheet "3_3" raw-data column (G) — shared-string text replaced in place
$ws1.Range("G2").Value = "6.893048,6.897519,6.900872,6.898537,6.897339,6.889532,6.892569,6.897221,6.894544,6.893426"
$ws1.Range("G3").Value = "6.892426,6.892029,6.899574,6.894904,6.892948,6.897041,6.897519,6.895302,6.903108,6.892408"
$ws1.Range("G4").Value = "6.895364,6.900413,6.895761,6.898537,6.900152,6.894165,6.899096,6.895085,6.902729,6.893066"

# --- Sheet "1_8": rows 2-4, columns B-E get the new smaller delta values ---
$ws2.Range("B2").Value = -0.001
$ws2.Range("C2").Value = 0.0065
$ws2.Range("D2").Value = -0.008800000000000001
$ws2.Range("E2").Value = 0.0047

$ws2.Range("B3").Value = -0.0009
$ws2.Range("C3").Value = 0.009299999999999999
$ws2.Range("D3").Value = -0.006
$ws2.Range("E3").Value = 0.0039

$ws2.Range("B4").Value = 0.0003
$ws2.Range("C4").Value = 0.0068
$ws2.Range("D4").Value = -0.004
$ws2.Range("E4").Value = 0.0036

# Sheet "1_8" raw-data column (G) — shared-string text replaced in place
$ws2.Range("G2").Value = "-0.002363,0.002897,0.006508,0.001557,0.002934,-0.008772,-0.002925,0.00152,-0.004333,-0.007377"
$ws2.Range("G3").Value = "0.00152,-0.001524,0.009305,-0.004043,-0.000968,-0.001826,-0.005975,-0.000425,-0.003209,-0.001542"
$ws2.Range("G4").Value = "0.001798,0.000977,0.002372,-0.004043,-0.004024,0.001508,-0.003487,-0.003191,0.004027,0.006774"
